$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A29").Value = 54
$ws.Range("B29").Value = "math"
$ws.Range("C29").Value = "riya-morankar"
$ws.Range("D29").Value = "N/A"
$ws.Range("E29").Value = "edit1 to main"

# Force the date-looking string to stay as text (matching the other
# Date column cells, which are plain text, not real Excel dates).
$ws.Range("F29").NumberFormat = "@"
$ws.Range("F29").Value = "2025-06-20"
$ws.Range("F29").Style = "Normal"
